$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting the existing rows 20-25 down to 21-26
$ws.Rows(20).Insert()

# Populate the newly inserted row 20 with the new record
$ws.Cells.Item(20, 1).Value = 5
$ws.Cells.Item(20, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(20, 3).Value = "Maule"
$ws.Cells.Item(20, 4).Value = 44466
$ws.Cells.Item(20, 5).Value = 7
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100107
$ws.Cells.Item(20, 8).Value = "Otros"
$ws.Cells.Item(20, 9).Value = 100107002
$ws.Cells.Item(20, 10).Value = "Chirimoya"
$ws.Cells.Item(20, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(20, 12).Value = "Especial"
$ws.Cells.Item(20, 13).Value = 110
$ws.Cells.Item(20, 14).Value = 30000
$ws.Cells.Item(20, 15).Value = 30000
$ws.Cells.Item(20, 16).Value = 30000
$ws.Cells.Item(20, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(20, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 19).Value = 3000
$ws.Cells.Item(20, 20).Value = 10

# Copy the date-column number format (style) from the row below onto the new D20 cell
$ws.Cells.Item(21, 4).Copy()
$ws.Cells.Item(20, 4).PasteSpecial(-4122)
